# Update countries & provincias Spain
#
# The source data table (rows 19-30 of the "provincias" sheet) got a new
# "Murcia" entry inserted right after "Zaragoza" (row 19), pushing the
# following provinces down by one row; the old "Murcia" row (which used
# to be the last row of this block, row 30) is dropped since its figures
# are superseded by the fresh ones now at row 19. A couple of individual
# counters elsewhere in the table were also refreshed, and the "last
# updated" timestamp moved from 22:16 to 22:46.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- timestamp banner (A1) ---------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 21 de Marzo de 2020 a las 22:46"

# --- single-cell tweaks --------------------------------------------------
$ws.Range("E13").Value = 8
$ws.Range("E45").Value = 2

# --- rows 19-30: insert the new Murcia figures and re-flow the rest ------
$data = @(
    @("Murcia",             296, 1,  213, 1),
    @("Granada",            289, 72, 276, 13),
    @("Gipuzkoa/Guipuzcoa", 273, 21, 262, 11),
    @("Burgos",             269, 27, 175, 14),
    @("Salamanca",          265, 13, 180, 15),
    @("Pontevedra",         264, 5,  262, 2),
    @("Guadalajara",        263, 2,  257, 4),
    @("Illes Balears",      246, 10, 232, 4),
    @("Sevilla",            245, 72, 243, 2),
    @("Caceres",            243, 2,  231, 10),
    @("Valladolid",         241, 13, 193, 10),
    @("Tenerife",           219, 4,  211, 5)
)

$startRow = 19
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
}
